$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.626.28"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.654.21"
$ws.Range("E3").Value = "  -2.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.96"
$ws.Range("E5").Value = "  +2.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3634"
$ws.Range("E7").Value = "  -2.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.81"
$ws.Range("E8").Value = "  -4.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3263"
$ws.Range("E9").Value = "  -4.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.134"
$ws.Range("E10").Value = "  -6.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07047"
$ws.Range("E11").Value = "  -5.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.019"
$ws.Range("E13").Value = "  -4.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.53"
$ws.Range("E14").Value = "  -6.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.651.63"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.612"
$ws.Range("E16").Value = "  -5.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001043"
$ws.Range("E17").Value = "  -7.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06636"
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.80"
$ws.Range("E20").Value = "  -5.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.902"
$ws.Range("E21").Value = "  -6.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.74"
$ws.Range("E22").Value = "  -8.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.53"
$ws.Range("E23").Value = "  -3.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.568.97"
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.440"
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.363"
$ws.Range("E26").Value = "  -14.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "148.00"
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.58"
$ws.Range("E28").Value = "  -8.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.836.42"
$ws.Range("E29").Value = "  -3.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.215"
$ws.Range("E30").Value = "  -2.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.53"
$ws.Range("E31").Value = "  -4.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.072"
$ws.Range("E32").Value = "  -3.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.812"
$ws.Range("E33").Value = "  -13.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08410"
$ws.Range("E34").Value = "  -3.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.676"
$ws.Range("E35").Value = "  -5.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.26"
$ws.Range("E36").Value = "  -10.37%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.213"
$ws.Range("E37").Value = "  -6.39%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.270"
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06027"
$ws.Range("E39").Value = "  -8.96%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02229"
$ws.Range("E40").Value = "  -7.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2070"
$ws.Range("E41").Value = "  -6.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.193"
$ws.Range("E42").Value = "  -9.45%  "
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5892"
$ws.Range("E44").Value = "  -7.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.822"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.75"
$ws.Range("E46").Value = "  -7.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5626"
$ws.Range("E47").Value = "  -7.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.81"
$ws.Range("E48").Value = "  -3.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.953"
$ws.Range("E49").Value = "  -6.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06953"
$ws.Range("E50").Value = "  -4.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.190"
$ws.Range("E51").Value = "  -2.42%  "
